# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets list the same set of events; "全部类型" simply has one extra
# row (row 17, inserted earlier) so its row numbers are offset by +1 from
# row 20 onward relative to "展览".

$wb = $excel.ActiveWorkbook

# Row -> (old, new) for the "展览" sheet (sheet1).
$sheet1Updates = @(
    @{Row=3;  Old=273;  New=275},
    @{Row=4;  Old=82;   New=83},
    @{Row=5;  Old=286;  New=287},
    @{Row=7;  Old=109;  New=110},
    @{Row=12; Old=119;  New=122},
    @{Row=13; Old=2498; New=2510},
    @{Row=14; Old=62;   New=67},
    @{Row=19; Old=539;  New=540},
    @{Row=20; Old=593;  New=597},
    @{Row=21; Old=177;  New=180},
    @{Row=22; Old=92;   New=93},
    @{Row=23; Old=50;   New=51},
    @{Row=24; Old=52;   New=53},
    @{Row=25; Old=2095; New=2108},
    @{Row=26; Old=4191; New=4206},
    @{Row=30; Old=1229; New=1230},
    @{Row=31; Old=243;  New=244},
    @{Row=32; Old=2130; New=2133},
    @{Row=34; Old=472;  New=473},
    @{Row=35; Old=67;   New=68},
    @{Row=36; Old=126;  New=129},
    @{Row=37; Old=295;  New=296},
    @{Row=38; Old=438;  New=439},
    @{Row=39; Old=725;  New=727},
    @{Row=42; Old=7;    New=11},
    @{Row=43; Old=432;  New=433}
)

# Row -> (old, new) for the "全部类型" sheet (sheet4).
$sheet4Updates = @(
    @{Row=3;  Old=273;  New=275},
    @{Row=4;  Old=82;   New=83},
    @{Row=5;  Old=286;  New=287},
    @{Row=7;  Old=109;  New=110},
    @{Row=12; Old=119;  New=122},
    @{Row=13; Old=2498; New=2510},
    @{Row=14; Old=62;   New=67},
    @{Row=20; Old=539;  New=540},
    @{Row=21; Old=593;  New=597},
    @{Row=22; Old=177;  New=180},
    @{Row=23; Old=92;   New=93},
    @{Row=24; Old=50;   New=51},
    @{Row=25; Old=52;   New=53},
    @{Row=26; Old=2095; New=2108},
    @{Row=27; Old=4191; New=4206},
    @{Row=31; Old=1229; New=1230},
    @{Row=32; Old=243;  New=244},
    @{Row=33; Old=2130; New=2133},
    @{Row=35; Old=472;  New=473},
    @{Row=36; Old=67;   New=68},
    @{Row=37; Old=126;  New=129},
    @{Row=38; Old=295;  New=296},
    @{Row=39; Old=438;  New=439},
    @{Row=40; Old=725;  New=727},
    @{Row=43; Old=7;    New=11},
    @{Row=44; Old=432;  New=433}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $sheet1Updates) {
    $cell = $ws1.Cells.Item($u.Row, 6)   # column F = "想去人数"
    $current = $cell.Value2
    if ($current -eq $u.Old) {
        $cell.Value = $u.New
    } else {
        $delta = $u.New - $u.Old
        $cell.Value = $current + $delta
    }
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $sheet4Updates) {
    $cell = $ws4.Cells.Item($u.Row, 6)   # column F = "想去人数"
    $current = $cell.Value2
    if ($current -eq $u.Old) {
        $cell.Value = $u.New
    } else {
        $delta = $u.New - $u.Old
        $cell.Value = $current + $delta
    }
}
